# Fruta / hortaliza, semanal
# Insert two new weekly price rows (week of 2022-01-17, serial 44578) for
# "Terminal La Palmera de La Serena" / Damasco / Modesto, ahead of the
# existing rows for 2022-01-10 (serial 44571), pushing the latter down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 58 - this shifts the former rows 58-59
# (date 44571) down to rows 60-61, unchanged.
$ws.Rows.Item(58).Insert()
$ws.Rows.Item(58).Insert()

# New row 58: Modesto / Especial, week of 2022-01-17
$ws.Range("A58").Value = 8
$ws.Range("B58").Value = "Terminal La Palmera de La Serena"
$ws.Range("C58").Value = "Coquimbo"
$ws.Range("D58").Value = 44578
$ws.Range("E58").Value = 4
$ws.Range("F58").Value = "Fruta"
$ws.Range("G58").Value = 100103
$ws.Range("H58").Value = "Frutos de hueso (carozo)"
$ws.Range("I58").Value = 100103003
$ws.Range("J58").Value = "Damasco"
$ws.Range("K58").Value = "Modesto"
$ws.Range("L58").Value = "Especial"
$ws.Range("M58").Value = 240
$ws.Range("N58").Value = 22500
$ws.Range("O58").Value = 23000
$ws.Range("P58").Value = 22750
$ws.Range("Q58").Value = "$/caja 18 kilos"
$ws.Range("R58").Value = "Región Metropolitana"
$ws.Range("S58").Value = 1264
$ws.Range("T58").Value = 18

# New row 59: Modesto / Primera, week of 2022-01-17
$ws.Range("A59").Value = 8
$ws.Range("B59").Value = "Terminal La Palmera de La Serena"
$ws.Range("C59").Value = "Coquimbo"
$ws.Range("D59").Value = 44578
$ws.Range("E59").Value = 4
$ws.Range("F59").Value = "Fruta"
$ws.Range("G59").Value = 100103
$ws.Range("H59").Value = "Frutos de hueso (carozo)"
$ws.Range("I59").Value = 100103003
$ws.Range("J59").Value = "Damasco"
$ws.Range("K59").Value = "Modesto"
$ws.Range("L59").Value = "Primera"
$ws.Range("M59").Value = 360
$ws.Range("N59").Value = 19500
$ws.Range("O59").Value = 20000
$ws.Range("P59").Value = 19750
$ws.Range("Q59").Value = "$/caja 18 kilos"
$ws.Range("R59").Value = "Región Metropolitana"
$ws.Range("S59").Value = 1097
$ws.Range("T59").Value = 18
